# Update sheet name and header label from "June 19" to "June 20"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-06-20"
$ws.Range("B1").Value = "June 2022 (through June 20)"

# New carjacking data for 2022-06-28 added to column B (and some corrections
# scattered through older months in the same rows)
$ws.Range("B2").Value = 7
$ws.Range("N2").Value = 6
$ws.Range("AF2").Value = 1
$ws.Range("AR2").Value = 2

$ws.Range("N5").Value = 2
$ws.Range("AF5").Value = 6

$ws.Range("AF8").Value = 1

$ws.Range("N9").Value = 4
$ws.Range("AL9").Value = 3

$ws.Range("N11").Value = 1

$ws.Range("N13").Value = 1

$ws.Range("N14").Value = 6

$ws.Range("H19").Value = 2

$ws.Range("H23").Value = 3

$ws.Range("H25").Value = 1

$ws.Range("N37").Value = 2

$ws.Range("AF39").Value = 1

$ws.Range("AF41").Value = 2

$ws.Range("T45").Value = 1

$ws.Range("H51").Value = 1
$ws.Range("N51").Value = 3

$ws.Range("H70").Value = 2
$ws.Range("AF70").Value = 2

$ws.Range("AF75").Value = 2

$ws.Range("Z89").Value = 1

$ws.Range("N95").Value = 1
